$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Total" row correct/total marks values
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 40
$ws.Range("E12").Value = "40/140"
